# Insert a new row at 212 (pushes the former rows 212-282 down to 213-283)
# and populate it with a new weekly price observation for Coliflor.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("212:212").Insert()

$ws.Cells.Item(212, 1).Value = 4
$ws.Cells.Item(212, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(212, 3).Value = "Los Lagos"
$ws.Cells.Item(212, 4).Value = 44627
$ws.Cells.Item(212, 5).Value = 10
$ws.Cells.Item(212, 6).Value = 100112008
$ws.Cells.Item(212, 7).Value = "Coliflor"
$ws.Cells.Item(212, 8).Value = "Sin especificar"
$ws.Cells.Item(212, 9).Value = "Segunda"
$ws.Cells.Item(212, 10).Value = 250
$ws.Cells.Item(212, 11).Value = 1800
$ws.Cells.Item(212, 12).Value = 1800
$ws.Cells.Item(212, 13).Value = 1800
$ws.Cells.Item(212, 14).Value = "$/unidad"
$ws.Cells.Item(212, 15).Value = "Región Metropolitana"
$ws.Cells.Item(212, 16).Value = 1800
$ws.Cells.Item(212, 17).Value = 1
$ws.Cells.Item(212, 18).Value = "Hortaliza"
